$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows -------------------------------------------------

# Row 3: "Bug in Simulation suchen" task is now done
$ws.Range("F3").Value = "done"

# Row 5: clarify the "Simulation fixen" task text and mark it "In Arbeit"
$ws.Range("B5").Value = "Simulation fixen (Hindernisse eintragen noch verbuggt)"
$ws.Range("F5").Value = "In Arbeit"

# Row 10: "Testbilder mit Tiefenkamera" got an end date and is now "offen"
$ws.Range("D5").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = 44274
$ws.Range("F10").Value = "offen"

# --- New rows --------------------------------------------------------------

# Row 12: Explodierende Voegel
$ws.Range("B12").Value = "Explodierende Vögel"
$ws.Range("C12").Value = "Martin"
$ws.Range("D5").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = 44273
$ws.Range("F12").Value = "offen"

# Row 13: Dronekit studieren
$ws.Range("B13").Value = "Dronekit studieren"
$ws.Range("D5").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = 44273
$ws.Range("F13").Value = "offen"

# Row 14: GPS der Drohne auslesen
$ws.Range("B14").Value = "GPS der Drohne auslesen"
$ws.Range("D5").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 44273
$ws.Range("F14").Value = "offen"

# --- View state: reflect the user's last selection / zoom ------------------
[void]$ws.Range("D15").Select()
$excel.ActiveWindow.Zoom = 86
